# Insert a new price record as row 29 on the single data sheet, pushing the
# existing rows 29-78 down to 30-79 (dimension grows from A1:R78 to A1:R79).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 29..78 down one row, leaving a blank row 29 to populate.
$ws.Rows("29:29").Insert()

# Populate the newly inserted row 29 with the new observation.
$ws.Range("A29").Value = 11
$ws.Range("B29").Value = "Vega Monumental Concepción"
$ws.Range("C29").Value = "Bíobío"
$ws.Range("D29").Value = 45259
$ws.Range("E29").Value = 8
$ws.Range("F29").Value = 100112026
$ws.Range("G29").Value = "Haba"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 150
$ws.Range("K29").Value = 8500
$ws.Range("L29").Value = 8500
$ws.Range("M29").Value = 8500
$ws.Range("N29").Value = "`$/saco 25 kilos"
$ws.Range("O29").Value = "Región del Maule"
$ws.Range("P29").Value = 340
$ws.Range("Q29").Value = 25
$ws.Range("R29").Value = "Hortaliza"
